# Update vm_pu results: slack bus voltage setpoint changed from 1.05 to 1.02 p.u.
# (reflecting the "case with 380 kV" scenario), recalculating downstream bus voltages.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @{
    2 = @(1.02,1.032015776714196,1.035866392573799,1.045687926806673,1.053467506576291,1.035366852227118,1.037148321012233,1.038661770439816,1.048455465856699,1.056213394052302,1.016353281504802)
    3 = @(1.02,1.032984689000636,1.036593769848524,1.046644822169625,1.054502041919094,1.035567820689132,1.037759252080038,1.039198849257479,1.049223473795543,1.057060412660611,1.016556932610429)
    4 = @(1.02,1.033611848989776,1.037064176104779,1.047264538805189,1.055171996345113,1.035696099101795,1.038154176128475,1.039545452585946,1.049720346485992,1.057608424677989,1.016688533721792)
    5 = @(1.02,1.033875556480089,1.037261872561194,1.047525196342901,1.055453773783258,1.035749605028806,1.038320108311775,1.039690942894676,1.049929211851236,1.057838792367596,1.016743816705194)
    6 = @(1.02,1.033919837009535,1.03729506297668,1.047568969408113,1.055501092970104,1.035758564127485,1.038347963529234,1.039715358330136,1.049964280073794,1.057877471120566,1.016753096482955)
    7 = @(1.02,1.033615372468678,1.037066817980917,1.047268021219581,1.055175760964853,1.035696815710767,1.038156393692236,1.039547397506914,1.049723137434492,1.057611502928475,1.016689272581597)
    8 = @(1.02,1.032343181537331,1.036112265325494,1.046011201422633,1.053817020111091,1.035435134866173,1.03735486818808,1.038843469249792,1.048715033773132,1.056499661046137,1.016422142367222)
    9 = @(1.02,1.030103044668473,1.034428316222941,1.043800711433209,1.051426926570699,1.034960553286086,1.035939528047279,1.037596029573687,1.046938050521128,1.054539987557991,1.015950098649647)
    10 = @(1.02,1.028610745870965,1.033304470099812,1.042329919969096,1.049836392600357,1.034635144762088,1.034994025078907,1.036759727210831,1.045753053894744,1.053233273009624,1.015634527455861)
    11 = @(1.02,1.027964837150423,1.032817558321574,1.041693742064922,1.049148362528585,1.034492106369851,1.034584159081206,1.036396499951994,1.045239865520509,1.052667397662469,1.015497677292472)
    12 = @(1.02,1.027724958336528,1.032636656749012,1.041457540897993,1.048892900656784,1.034438655365709,1.034431848534576,1.036261416230448,1.045049233610237,1.052457198139214,1.01544681438099)
    13 = @(1.02,1.027776411319722,1.032675462582433,1.041508202147799,1.048947693416955,1.034450135268429,1.034464522734542,1.036290399599527,1.045090125284508,1.052502287011589,1.015457726032014)
    14 = @(1.02,1.027945007865449,1.032802605768561,1.041674215484129,1.049127243862541,1.034487694621877,1.034571570429029,1.036385337249989,1.04522410804578,1.052650022668088,1.015493473572538)
    15 = @(1.02,1.028048891132995,1.03288093749917,1.041776515535781,1.049237884505513,1.034510793767151,1.034637517044833,1.036443809635062,1.045306657822927,1.052741046412808,1.015515494755629)
    16 = @(1.02,1.028653618354985,1.033336779043256,1.042372155490896,1.049882069313283,1.034644592822267,1.03502121695776,1.036783810235989,1.045787110951806,1.053270827101713,1.015643605436405)
    17 = @(1.02,1.029033019401465,1.033622642454967,1.042745968520485,1.050286332428173,1.034727950351608,1.03526177984803,1.036996788899292,1.046088466362362,1.053603129135598,1.015723910968678)
    18 = @(1.02,1.0292543432998,1.033789354757509,1.04296407332802,1.050522197982801,1.034776365421543,1.035402051981272,1.037120909295003,1.046264234399001,1.053796949519793,1.015770731958359)
    19 = @(1.02,1.02932981337425,1.033846194734383,1.043038452600239,1.050602633166707,1.034792838758885,1.035449873650464,1.037163213063113,1.046324165468727,1.053863036262335,1.015786693328354)
    20 = @(1.02,1.028992310599868,1.033591974800186,1.042705855088495,1.050242952018883,1.034719028178141,1.035235974286892,1.036973949310086,1.04605613454714,1.053567476864005,1.015715296998388)
    21 = @(1.02,1.027895359265188,1.03276516638864,1.041625325827228,1.049074367878757,1.034476643162978,1.034540049425325,1.036357385025545,1.045184653757477,1.052606518420953,1.015482947657452)
    22 = @(1.02,1.027205896620537,1.032245083249393,1.040946554171653,1.048340229692625,1.034322393879862,1.034102100326145,1.035968772631807,1.044636656472889,1.052002278467401,1.015336683092389)
    23 = @(1.02,1.027571371479126,1.032520811168687,1.041306326733805,1.048729353401353,1.034404339749339,1.034334302534636,1.036174873592328,1.044927165923385,1.052322601747604,1.015414237431475)
    24 = @(1.02,1.029010705075593,1.033605832276763,1.042723980417857,1.050262553555479,1.034723060360694,1.035247634845541,1.036984269866281,1.046070743925832,1.053583586608724,1.015719189342176)
    25 = @(1.02,1.030681977695436,1.034863876228348,1.04437167438088,1.05204432232094,1.035084836403859,1.036305772791401,1.037919350392011,1.047397506769714,1.055046660421646,1.01607228874276)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}

Write-Host "Updated vm_pu values for rows 2-25"
